$wb = $excel.ActiveWorkbook

# --- Sheet: ip_address_list ---
$ws1 = $wb.Worksheets.Item("ip_address_list")

# D5 (514_Teleflex): merge "CAM: 192.168.14.???" and "NAS:192.168.14.245" lines into one,
# also dropping one of the "?" characters.
$ws1.Range("D5").Value = "PC:192.168.14.240`nCAM: 192.168.14.??NAS:192.168.14.245`n*******************************`nuser: Vision`npass: *Jhv2708"

# D12 (515_ZF Stara Boleslav): join first two lines, drop an "s" from "User" and "ad" from "jhvadm1n".
$ws1.Range("D12").Value = "NAS - 10.9.250.100Uer:spravce Pass:Jhv*2708 `nUser:jhvadmin Pass:jhvm1n >>> na portu 8080. `n123TPV456"

# --- Sheet: ip_adress_fav_list ---
$ws2 = $wb.Worksheets.Item("ip_adress_fav_list")

# D2 (514_Teleflex): same edit as ws1 D5.
$ws2.Range("D2").Value = "PC:192.168.14.240`nCAM: 192.168.14.??NAS:192.168.14.245`n*******************************`nuser: Vision`npass: *Jhv2708"

# D5 (515_ZF Stara Boleslav): same edit as ws1 D12.
$ws2.Range("D5").Value = "NAS - 10.9.250.100Uer:spravce Pass:Jhv*2708 `nUser:jhvadmin Pass:jhvm1n >>> na portu 8080. `n123TPV456"

# --- Sheet: Settings ---
$ws4 = $wb.Worksheets.Item("Settings")

# B3 ("spousteci okno: na oblibenych = 1"): flip flag from 1 to 0.
$ws4.Range("B3").Value = 0
